# Apply the "Integration plan" update:
#  - Row 8 gets a "Done" marker in column M (status column)
#  - Row 9 gets fully filled in with T/S/X test status values
#  - Rows 10-12 get their remaining status columns filled in (replacing the
#    Danish planning notes that used to live in columns E/F) with T/S/X values
#  - Row 13 gets fully filled in with X values
#  - Rows 14-15 lose the now-obsolete Danish planning notes
#  - Selection moves to L13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 8 - add the "Done" status in column M
$ws.Range("M8").Value = "Done"

# Row 9 - fill out the whole row with statuses
$ws.Range("B9").Value = "T"
$ws.Range("C9").Value = "X"
$ws.Range("D9").Value = "S"
$ws.Range("E9").Value = "S"
$ws.Range("F9").Value = "S"
$ws.Range("G9").Value = "X"
$ws.Range("H9").Value = "X"
$ws.Range("I9").Value = "X"
$ws.Range("J9").Value = "X"
$ws.Range("K9").Value = "X"
$ws.Range("L9").Value = "S"

# Row 10 - replace old notes in E10/F10, fill remaining new columns
$ws.Range("B10").Value = "T"
$ws.Range("D10").Value = "S"
$ws.Range("E10").Value = "S"
$ws.Range("F10").Value = "X"
$ws.Range("G10").Value = "X"

# Row 11 - replace old notes in E11/F11, fill remaining new columns
$ws.Range("B11").Value = "T"
$ws.Range("D11").Value = "X"
$ws.Range("E11").Value = "S"
$ws.Range("F11").Value = "X"
$ws.Range("G11").Value = "X"

# Row 12 - replace old notes in E12/F12, fill remaining new columns
$ws.Range("B12").Value = "T"
$ws.Range("D12").Value = "X"
$ws.Range("E12").Value = "X"
$ws.Range("F12").Value = "X"
$ws.Range("G12").Value = "X"
$ws.Range("L12").Value = "S"

# Row 13 - fully filled in with X, replacing old notes in E13/F13
$ws.Range("B13").Value = "X"
$ws.Range("C13").Value = "X"
$ws.Range("D13").Value = "X"
$ws.Range("E13").Value = "X"
$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"
$ws.Range("H13").Value = "X"
$ws.Range("I13").Value = "X"
$ws.Range("J13").Value = "X"
$ws.Range("K13").Value = "X"
$ws.Range("L13").Value = "X"

# Row 14/15 - clear the now-obsolete Danish planning notes
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("E15").Value = ""

# Update the active selection to match the author's final cursor position
$ws.Range("L13").Select()
